# Update column G ("K") values for rows 2-20 as computed by the
# regenerated save_data (uses K count instead of Strike# values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    20 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
